$d = $word.ActiveDocument

# Start from the last paragraph in the document body (the one ending
# "...pH-Wert zwischen 3 und 8. ") and append six new paragraphs after it,
# mirroring the existing "Bodenreaktion" heading/body pattern.

$texts = @(
    "pH-Wert im Boden messen",
    "1. man gräbt ein kleines Loch in den Boden",
    "2.  man füllt das Loch mit Wasser auf",
    "3. man steckt die Prüfspitze des Geräts für 60 Sekunden in den Schlamm",
    "4. der pH-Wert kann am Gerät abgelesen werden",
    ""
)
$bolds = @($true, $false, $false, $false, $false, $true)

for ($i = 0; $i -lt $texts.Length; $i++) {
    $p = $d.Paragraphs.Last
    $p.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    if ($texts[$i] -ne "") {
        $newPara.Range.Text = $texts[$i]
    }
    $newPara.Range.Font.Bold = $bolds[$i]
}
